$wb = $excel.ActiveWorkbook

# --- Window view change (cosmetic, workbook-level) ---
$wb.Windows.Item(1).Left = 4080
$wb.Windows.Item(1).Top = 0
$wb.Windows.Item(1).Width = 24720
$wb.Windows.Item(1).Height = 16740

# --- Sheet 1: openbis-metadata ---
$ws1 = $wb.Worksheets.Item("openbis-metadata")

# Row 3: Strain / strain1 / description -> change Value to "MGP9"
$ws1.Range("B3").Value = "MGP9"

# Row 2: Experiment / (empty) / description -> set Value to "/TEST/TEST/TEST"
$ws1.Range("B2").Value = "/TEST/TEST/TEST"

# Update the selected cell on the sheet to B3
$ws1.Range("B3").Select()

# --- Sheet 2: openbis-data (no content changes besides shared-string reindex, already handled) ---
